$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.73"
$ws.Range("E2").Value = "'1.70%"
$ws.Range("D3").Value = "'41.72"
$ws.Range("E3").Value = "'3.25%"
$ws.Range("D4").Value = "'5.012"
$ws.Range("E4").Value = "'-0.19%"
$ws.Range("D5").Value = "'0.07521"
$ws.Range("E5").Value = "'2.99%"
$ws.Range("D6").Value = "'1.605"
$ws.Range("E6").Value = "'4.49%"
$ws.Range("D7").Value = "'0.9173"
$ws.Range("E7").Value = "'-1.04%"
$ws.Range("D9").Value = "'0.1176"
$ws.Range("E9").Value = "'0.99%"
$ws.Range("D10").Value = "'0.1830"
$ws.Range("E10").Value = "'4.37%"
$ws.Range("D11").Value = "'0.08961"
$ws.Range("E11").Value = "'3.49%"
$ws.Range("D12").Value = "'0.04098"
$ws.Range("E12").Value = "'-5.86%"
$ws.Range("D13").Value = "'0.1049"
$ws.Range("E13").Value = "'-0.44%"
$ws.Range("D14").Value = "'0.001277"
$ws.Range("E14").Value = "'0.44%"
$ws.Range("D15").Value = "'0.005983"
$ws.Range("E15").Value = "'-0.45%"
$ws.Range("E16").Value = "'-0.04%"
$ws.Range("E17").Value = "'2.11%"
$ws.Range("D18").Value = "'0.3329"
$ws.Range("E18").Value = "'1.47%"
$ws.Range("D19").Value = "'8.300"
$ws.Range("E19").Value = "'4.05%"
$ws.Range("E20").Value = "'-2.78%"
$ws.Range("D21").Value = "'0.3103"
$ws.Range("E21").Value = "'11.80%"
$ws.Range("D22").Value = "'0.04100"
$ws.Range("E22").Value = "'4.05%"
$ws.Range("D23").Value = "'0.001266"
$ws.Range("E23").Value = "'0.20%"
$ws.Range("D24").Value = "'0.003893"
$ws.Range("E24").Value = "'6.32%"
$ws.Range("E25").Value = "'8.23%"
$ws.Range("D38").Value = "'0.02398"
$ws.Range("D39").Value = "'0.05205"
$ws.Range("E39").Value = "'2.85%"
$ws.Range("D40").Value = "'0.006305"
$ws.Range("E40").Value = "'7.68%"
$ws.Range("D41").Value = "'0.007802"
$ws.Range("E41").Value = "'-0.81%"
$ws.Range("D42").Value = "'0.1326"
$ws.Range("E42").Value = "'3.01%"
$ws.Range("D43").Value = "'0.007411"
$ws.Range("E43").Value = "'0.76%"
$ws.Range("D44").Value = "'0.007590"
$ws.Range("E44").Value = "'4.37%"
$ws.Range("D45").Value = "'0.3241"
$ws.Range("E45").Value = "'1.79%"
$ws.Range("D46").Value = "'0.00006583"
$ws.Range("E46").Value = "'6.19%"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.04531"
$ws.Range("E48").Value = "'-15.94%"
$ws.Range("D49").Value = "'0.004203"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.11%"
